$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-of-the-day note in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText.Replace("1000 Bs = 5.13 = 20333.44 pesos", "1000 Bs = 5.1 = 20236.58 pesos")
$newText = $newText.Replace("20333.44 pesos = 5.11 = 953.16 Bs", "20236.58 pesos = 5.07 = 949.22 Bs")
$wsHoja1.Range("A1").Value2 = $newText

# --- Sheet "tasas": update the rate table values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 196.13
$wsTasas.Range("O10").Value = 3969
$wsTasas.Range("N12").Value = 3990
$wsTasas.Range("O12").Value = 187.156
